$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.722.45'
$ws.Range("E2").Value = '  +3.79%  '

$ws.Range("D3").Value = '2.739.90'
$ws.Range("E3").Value = '  +3.25%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '579.20'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").Value = '159.20'
$ws.Range("E6").Value = '  +9.90%  '

$ws.Range("E7").Value = '  +3.12%  '

$ws.Range("D8").Value = '0.996'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '2.759.65'
$ws.Range("E9").Value = '  +3.34%  '

$ws.Range("E10").Value = '  +3.50%  '

$ws.Range("E11").Value = '  +2.75%  '

$ws.Range("E12").Value = '  +3.67%  '

$ws.Range("E13").Value = '  +0.47%  '

$ws.Range("D14").Value = '3.230.25'
$ws.Range("E14").Value = '  +3.65%  '

$ws.Range("E15").Value = '  +5.17%  '

$ws.Range("D16").Value = '63.706.02'
$ws.Range("E16").Value = '  +3.94%  '

$ws.Range("D17").Value = '0.0000156'
$ws.Range("E17").Value = '  +6.66%  '

$ws.Range("D18").Value = '2.753.71'
$ws.Range("E18").Value = '  +3.72%  '

$ws.Range("D19").Value = '12.16'
$ws.Range("E19").Value = '  +3.38%  '

$ws.Range("D20").Value = '4.94'
$ws.Range("E20").Value = '  +3.49%  '

$ws.Range("D21").Value = '364.22'
$ws.Range("E21").Value = '  +2.76%  '

$ws.Range("D22").Value = '7.03'
$ws.Range("E22").Value = '  +1.93%  '

$ws.Range("D23").Value = '0.540'
$ws.Range("E23").Value = '  +2.47%  '

$ws.Range("D24").Value = '0.994'
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("D25").Value = '66.64'
$ws.Range("E25").Value = '  +3.51%  '

$ws.Range("E26").Value = '  +5.51%  '

$ws.Range("D27").Value = '8.63'
$ws.Range("E27").Value = '  +1.02%  '

$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").Value = '0.0₃0918'
$ws.Range("E29").Value = '  +11.55%  '

$ws.Range("E30").Value = '  +1.23%  '

$ws.Range("D31").Value = '7.27'
$ws.Range("E31").Value = '  +5.51%  '

$ws.Range("E32").Value = '  +14.01%  '

$ws.Range("D33").Value = '173.75'
$ws.Range("E33").Value = '  +2.95%  '

$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = '20.66'
$ws.Range("E35").Value = '  +2.89%  '

$ws.Range("E36").Value = '  +5.82%  '

$ws.Range("E37").Value = '  +6.49%  '

$ws.Range("E38").Value = '  +6.55%  '

$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("D40").Value = '4.26'
$ws.Range("E40").Value = '  +2.78%  '

$ws.Range("D41").Value = '338.78'
$ws.Range("E41").Value = '  -0.34%  '

$ws.Range("D42").Value = '6.23'
$ws.Range("E42").Value = '  +16.08%  '

$ws.Range("D43").Value = '39.63'
$ws.Range("E43").Value = '  +3.16%  '

$ws.Range("D44").Value = '22.51'
$ws.Range("E44").Value = '  +6.56%  '

$ws.Range("D45").Value = '21.89'
$ws.Range("E45").Value = '  +5.83%  '

$ws.Range("E46").Value = '  +4.01%  '

$ws.Range("D47").Value = '0.646'
$ws.Range("E47").Value = '  +2.44%  '

$ws.Range("D48").Value = '0.0260'
$ws.Range("E48").Value = '  +3.21%  '

$ws.Range("D49").Value = '137.53'
$ws.Range("E49").Value = '  +1.46%  '

$ws.Range("E50").Value = '  +2.25%  '

$ws.Range("E51").Value = '  +0.02%  '
